$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Process edits from the bottom of the document upward so that paragraph
# indices used for earlier edits are not invalidated by insertions made for
# later (lower-in-document) edits.
# ---------------------------------------------------------------------------

# --- FIELD DIRECTOR (The Feldman Group) -------------------------------------
# Add a new bullet after "Enhanced value of research deliverables..."
$pFeldman = $d.Paragraphs.Item(69)
if ($pFeldman.Range.Text -notlike "*Enhanced value of research deliverables*") {
    throw "Unexpected paragraph 69: " + $pFeldman.Range.Text
}
$pFeldman.Range.InsertParagraphAfter()
$d.Paragraphs.Item(70).Range.Text = "• Trained staff on PHP/MySQL for data analysis and reporting systems"

# --- PROGRAMMER (Lake Research Partners) ------------------------------------
# Add a new bullet after "Developed innovative approaches to visualizing..."
$pLake = $d.Paragraphs.Item(61)
if ($pLake.Range.Text -notlike "*Developed innovative approaches to visualizing*") {
    throw "Unexpected paragraph 61: " + $pLake.Range.Text
}
$pLake.Range.InsertParagraphAfter()
$d.Paragraphs.Item(62).Range.Text = "• Trained staff on building Python tooling for report generation and analysis"

# --- INTERIM TECHNOLOGY MANAGER (The Praxis Project) ------------------------
# Replace the 4 existing bullets and append 3 additional new bullets.
if ($d.Paragraphs.Item(50).Range.Text -notlike "*Integrated technology solutions within organizational frameworks*") {
    throw "Unexpected paragraph 50: " + $d.Paragraphs.Item(50).Range.Text
}
$d.Paragraphs.Item(50).Range.Text = "• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO"
$d.Paragraphs.Item(51).Range.Text = "• Directed all technology decisions and practices for massive multinational non-governmental organization"
$d.Paragraphs.Item(52).Range.Text = "• Developed comprehensive frameworks for internal and external technology audits"
$d.Paragraphs.Item(53).Range.Text = "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research"

$d.Paragraphs.Item(53).Range.InsertParagraphAfter()
$d.Paragraphs.Item(54).Range.Text = "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL"

$d.Paragraphs.Item(54).Range.InsertParagraphAfter()
$d.Paragraphs.Item(55).Range.Text = "• Managed technology infrastructure supporting community health initiatives across multiple countries"

$d.Paragraphs.Item(55).Range.InsertParagraphAfter()
$d.Paragraphs.Item(56).Range.Text = "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# --- SOFTWARE ENGINEER (Salsa Labs, Inc.) -----------------------------------
# Replace the 4 existing bullets and append 2 additional new bullets.
if ($d.Paragraphs.Item(44).Range.Text -notlike "*Developed software solutions for political campaigns*") {
    throw "Unexpected paragraph 44: " + $d.Paragraphs.Item(44).Range.Text
}
$d.Paragraphs.Item(44).Range.Text = "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously"
$d.Paragraphs.Item(45).Range.Text = "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"
$d.Paragraphs.Item(46).Range.Text = "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"
$d.Paragraphs.Item(47).Range.Text = "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"

$d.Paragraphs.Item(47).Range.InsertParagraphAfter()
$d.Paragraphs.Item(48).Range.Text = "• Collaborated with political strategists to translate geospatial requirements into technical solutions"

$d.Paragraphs.Item(48).Range.InsertParagraphAfter()
$d.Paragraphs.Item(49).Range.Text = "• Handled billions of records with millions of columns in high-performance CRM system"

# --- RESEARCH DIRECTOR (Progressive Change Campaign Committee) -------------
$d.Content.Find.Execute(
    "Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2) | Out-Null

# --- PROFESSIONAL SUMMARY ----------------------------------------------------
$d.Content.Find.Execute(
    "21 years", $true, $false, $false, $false, $false, $true, 1, $false,
    "15+ years", 2) | Out-Null

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
